$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoices")

# Insert a new "Building" column at position F (shifts City..Client Project right by one)
$ws.Columns.Item(6).Insert()

# Set header for new Building column
$ws.Range("F1").Value = "Building"

# Populate Building (F) and fix City/Province/Postal/Agreement/Project (G:K) per row
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "Stoney Creek"
$ws.Range("H2").Value = "ON"
$ws.Range("I2").Value = "L8E 0J7"
$ws.Range("J2").Value = "A0332"
$ws.Range("K2").Value = "The Shores"

$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "Vaughan"
$ws.Range("H3").Value = "ON"
$ws.Range("I3").Value = "L4K 4B4"
$ws.Range("J3").Value = "A0212"
$ws.Range("K3").Value = "Origins - Additional Lots"

$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = "Burlington"
$ws.Range("H4").Value = "ON"
$ws.Range("I4").Value = "L7L 6A9"
$ws.Range("J4").Value = "A0224"
$ws.Range("K4").Value = "Casa De Torri"

$ws.Range("F5").Value = "Bldg. A"
$ws.Range("G5").Value = "Richmond Hill"
$ws.Range("H5").Value = "ON"
$ws.Range("I5").Value = "L4B 1B9"
$ws.Range("J5").Value = "A0178"
$ws.Range("K5").Value = "Urban North"

$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "Burlington"
$ws.Range("H6").Value = "ON"
$ws.Range("I6").Value = "L7L 6A9"
$ws.Range("J6").Value = "A0224"
$ws.Range("K6").Value = "Casa De Torri"

$ws.Range("F7").Value = "Bldg. A"
$ws.Range("G7").Value = "Richmond Hill"
$ws.Range("H7").Value = "ON"
$ws.Range("I7").Value = "L4B 1B9"
$ws.Range("J7").Value = "A0178"
$ws.Range("K7").Value = "Urban North"

$ws.Range("F8").Value = "N/A"
$ws.Range("G8").Value = "Burlington"
$ws.Range("H8").Value = "ON"
$ws.Range("I8").Value = "L7L 6A9"
$ws.Range("J8").Value = "A0224"
$ws.Range("K8").Value = "Casa De Torri"

$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "Burlington"
$ws.Range("H9").Value = "ON"
$ws.Range("I9").Value = "L7L 6A9"
$ws.Range("J9").Value = "A0224"
$ws.Range("K9").Value = "Casa De Torri"

$ws.Range("F10").Value = "Bldg. A"
$ws.Range("G10").Value = "Richmond Hill"
$ws.Range("H10").Value = "ON"
$ws.Range("I10").Value = "L4B 1B9"
$ws.Range("J10").Value = "A0178"
$ws.Range("K10").Value = "Urban North"

$ws.Range("F11").Value = "Bldg. A"
$ws.Range("G11").Value = "Richmond Hill"
$ws.Range("H11").Value = "ON"
$ws.Range("I11").Value = "L4B 1B9"
$ws.Range("J11").Value = "A0178"
$ws.Range("K11").Value = "Urban North"

$ws.Range("F12").Value = "N/A"
$ws.Range("G12").Value = "Stoney Creek"
$ws.Range("H12").Value = "ON"
$ws.Range("I12").Value = "L8E 0J7"
$ws.Range("J12").Value = "N/A"
$ws.Range("K12").Value = "N/A"

$ws.Range("F13").Value = "N/A"
$ws.Range("G13").Value = "Vaughan"
$ws.Range("H13").Value = "ON"
$ws.Range("I13").Value = "L4L 8A9"
$ws.Range("J13").Value = "A0451"
$ws.Range("K13").Value = "North West"

$ws.Range("F14").Value = "N/A"
$ws.Range("G14").Value = "Vaughan"
$ws.Range("H14").Value = "ON"
$ws.Range("I14").Value = "L4K 4B4"
$ws.Range("J14").Value = "A0429-DE"
$ws.Range("K14").Value = "GO Towns"

$ws.Range("F15").Value = "N/A"
$ws.Range("G15").Value = "Vaughan"
$ws.Range("H15").Value = "ON"
$ws.Range("I15").Value = "L4K 4B4"
$ws.Range("J15").Value = "A0456DE"
$ws.Range("K15").Value = "Whitby Meadows Phase 3"

$ws.Range("F16").Value = "N/A"
$ws.Range("G16").Value = "Concord"
$ws.Range("H16").Value = "ON"
$ws.Range("I16").Value = "L4K 5R2"
$ws.Range("J16").Value = "A0426"
$ws.Range("K16").Value = "Erin Glen"

$ws.Range("F17").Value = "N/A"
$ws.Range("G17").Value = "Barrie"
$ws.Range("H17").Value = "ON"
$ws.Range("I17").Value = "L4M 0J4"
$ws.Range("J17").Value = "A0330"
$ws.Range("K17").Value = "Windfall Phase 4B"

$ws.Range("F18").Value = "N/A"
$ws.Range("G18").Value = "Vaughan"
$ws.Range("H18").Value = "ON"
$ws.Range("I18").Value = "L4K 4B4"
$ws.Range("J18").Value = "A0504"
$ws.Range("K18").Value = "Greenwood Seaton"

$ws.Range("F19").Value = "N/A"
$ws.Range("G19").Value = "Concord"
$ws.Range("H19").Value = "ON"
$ws.Range("I19").Value = "L4K 5R2"
$ws.Range("J19").Value = "A0426"
$ws.Range("K19").Value = "Erin Glen"

$ws.Range("F20").Value = "N/A"
$ws.Range("G20").Value = "N/A"
$ws.Range("H20").Value = "N/A"
$ws.Range("I20").Value = "L4K 3Z9"
$ws.Range("J20").Value = "A0501"
$ws.Range("K20").Value = "Fairway Meadows Phase 2"

$ws.Range("F21").Value = "N/A"
$ws.Range("G21").Value = "Vaughan"
$ws.Range("H21").Value = "ON"
$ws.Range("I21").Value = "L4K 4B4"
$ws.Range("J21").Value = "A0546"
$ws.Range("K21").Value = "Honeystone"

$ws.Range("F22").Value = "N/A"
$ws.Range("G22").Value = "Vaughan"
$ws.Range("H22").Value = "ON"
$ws.Range("I22").Value = "L4K 4B4"
$ws.Range("J22").Value = "A0546"
$ws.Range("K22").Value = "Honeystone"

$ws.Range("F23").Value = "N/A"
$ws.Range("G23").Value = "Vaughan"
$ws.Range("H23").Value = "ON"
$ws.Range("I23").Value = "L4K 4B4"
$ws.Range("J23").Value = "A0546"
$ws.Range("K23").Value = "Honeystone"

$ws.Range("F24").Value = "N/A"
$ws.Range("G24").Value = "Mississauga"
$ws.Range("H24").Value = "ON"
$ws.Range("I24").Value = "L5N 6C3"
$ws.Range("J24").Value = "A0521"
$ws.Range("K24").Value = "Woodstock"

$ws.Range("F25").Value = "N/A"
$ws.Range("G25").Value = "Mississauga"
$ws.Range("H25").Value = "ON"
$ws.Range("I25").Value = "L5N 6C3"
$ws.Range("J25").Value = "A0522"
$ws.Range("K25").Value = "New Dundee Kitchener"

$ws.Range("F26").Value = "N/A"
$ws.Range("G26").Value = "Vaughan"
$ws.Range("H26").Value = "ON"
$ws.Range("I26").Value = "L4K 4B4"
$ws.Range("J26").Value = "A0429-DE"
$ws.Range("K26").Value = "GO Towns"

$ws.Range("F27").Value = "N/A"
$ws.Range("G27").Value = "Vaughan"
$ws.Range("H27").Value = "ON"
$ws.Range("I27").Value = "L4K 4B4"
$ws.Range("J27").Value = "A0445"
$ws.Range("K27").Value = "Seatonville"

$ws.Range("F28").Value = "N/A"
$ws.Range("G28").Value = "Vaughan"
$ws.Range("H28").Value = "ON"
$ws.Range("I28").Value = "L4K 4B4"
$ws.Range("J28").Value = "A0464"
$ws.Range("K28").Value = "Seaton South"

$ws.Range("F29").Value = "N/A"
$ws.Range("G29").Value = "Vaughan"
$ws.Range("H29").Value = "ON"
$ws.Range("I29").Value = "L4K 4B4"
$ws.Range("J29").Value = "A0504"
$ws.Range("K29").Value = "Greenwood Seaton"

$ws.Range("F30").Value = "N/A"
$ws.Range("G30").Value = "Vaughan"
$ws.Range("H30").Value = "ON"
$ws.Range("I30").Value = "L4K 4B4"
$ws.Range("J30").Value = "A0546"
$ws.Range("K30").Value = "Honeystone"
